$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy the header style (bold font, border, centered) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I and J columns
$data = @(
    @(2, 7, 7),
    @(3, 9, 9),
    @(4, 7, 7),
    @(5, 6, 6),
    @(6, 7, 7),
    @(7, 7, 7),
    @(8, 6, 7),
    @(9, 7, 7),
    @(10, 7, 7),
    @(11, 7, 7),
    @(12, 7, 7),
    @(13, 7, 7),
    @(14, 10, 10),
    @(15, 8, 8),
    @(16, 9, 9),
    @(17, 7, 8),
    @(18, 7, 7),
    @(19, 6, 6),
    @(20, 8, 8),
    @(21, 6, 7),
    @(22, 8, 8),
    @(23, 5, 6),
    @(24, 7, 7),
    @(25, 6, 7),
    @(26, 6, 6),
    @(27, 10, 10),
    @(28, 10, 10),
    @(29, 6, 6),
    @(30, 7, 7),
    @(31, 8, 8),
    @(32, 6, 7),
    @(33, 7, 7),
    @(34, 7, 7),
    @(35, 7, 7),
    @(36, 7, 7),
    @(37, 8, 8),
    @(38, 7, 7),
    @(39, 6, 6),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 6, 7),
    @(43, 7, 7),
    @(44, 9, 9),
    @(45, 7, 7),
    @(46, 6, 6),
    @(47, 8, 8),
    @(48, 7, 7),
    @(49, 8, 8),
    @(50, 7, 7),
    @(51, 9, 9),
    @(52, 7, 7),
    @(53, 6, 7),
    @(54, 7, 7),
    @(55, 8, 8),
    @(56, 9, 9),
    @(57, 10, 10),
    @(58, 8, 8),
    @(59, 6, 6),
    @(60, 7, 7),
    @(61, 5, 5),
    @(62, 6, 6),
    @(63, 7, 7),
    @(64, 8, 8),
    @(65, 5, 5),
    @(66, 7, 7),
    @(67, 7, 7)
)

foreach ($row in $data) {
    $r = $row[0]
    $i = $row[1]
    $j = $row[2]
    $ws.Cells.Item($r, 9).Value = $i
    $ws.Cells.Item($r, 10).Value = $j
}

Write-Host "Applied I0/IF columns"